$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9655929803848267
$ws.Range("B1").Value = 3.257440328598022
$ws.Range("C1").Value = 4.084120273590088
$ws.Range("D1").Value = 3.091893672943115
$ws.Range("E1").Value = 1.340919375419617
